$d = $word.ActiveDocument

# Replace the two most-specific "June 22, 2022" contexts first so the
# generic "June 22, 2022" search afterwards only matches the one
# remaining (unqualified) occurrence.

# 1. " on June 22, 2022." -> " on June 24, 2022."
$d.Content.Find.Execute(" on June 22, 2022.", $true, $false, $false, $false, $false,
                         $true, 1, $false, " on June 24, 2022.", 2)

# 2. " license is suspended from June 22, 2022" -> " license is suspended from June 24, 2022"
$d.Content.Find.Execute(" license is suspended from June 22, 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, " license is suspended from June 24, 2022", 2)

# 3. remaining bare "June 22, 2022" (bold, "in full by") -> "June 24, 2022"
$d.Content.Find.Execute("June 22, 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "June 24, 2022", 2)

# 4. "August 21, 2022" -> "August 23, 2022"
$d.Content.Find.Execute("August 21, 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "August 23, 2022", 2)
